$d = $word.ActiveDocument

# --- Edit 1 --------------------------------------------------------------
# "...create and simulate our react native app. ..."
#   -> "...create and simulate our React Native app. ..."
# (capitalise "React Native" as a proper noun / product name)
$found1 = $d.Content.Find.Execute(
    "simulate our react native app",   # FindText
    $true,                             # MatchCase
    $false,                            # MatchWholeWord
    $false,                            # MatchWildcards
    $false,                            # MatchSoundsLike
    $false,                            # MatchAllWordForms
    $true,                             # Forward
    1,                                 # Wrap (wdFindContinue)
    $false,                            # Format
    "simulate our React Native app",   # ReplaceWith
    2)                                 # Replace (wdReplaceAll)

if (-not $found1) {
    throw "Edit 1 (React Native capitalisation) target text not found"
}

# --- Edit 2 --------------------------------------------------------------
# "In order to help the user track their moles they must be able to take
#  periodic photos so that a doctor can look for changes."
#   -> "To help the user track their moles, they must be able to take
#       periodic photos so that a doctor can look for changes."
$found2 = $d.Content.Find.Execute(
    "In order to help the user track their moles they must",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "To help the user track their moles, they must",
    2)

if (-not $found2) {
    throw "Edit 2 (sentence reword) target text not found"
}
